$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 1302
$ws.Range("I61").Value = 1302
$ws.Range("K61").Value = 3906
$ws.Range("M61").Value = -3734
$ws.Range("H74").Value = 6969.0435
$ws.Range("I74").Value = 4660.923
$ws.Range("K74").Value = 4660.923
$ws.Range("M74").Value = -3724.923
$ws.Range("H77").Value = 6969.0435
$ws.Range("I77").Value = 4660.923
$ws.Range("K77").Value = 23304.615
$ws.Range("M77").Value = -18624.615
$ws.Range("H132").Value = 3401.1897
$ws.Range("I132").Value = 2473.283
$ws.Range("K132").Value = 7419.849
$ws.Range("M132").Value = -4889.849
$ws.Range("H137").Value = 71430776
$ws.Range("I137").Value = 125000610
$ws.Range("J137").Value = 4332.6665
$ws.Range("K137").Value = 375001830
$ws.Range("L137").Value = 12997.9995
$ws.Range("M137").Value = -374999280
$ws.Range("N137").Value = -18097.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3693.5386
$ws.Range("J61").Value = 5675
$ws.Range("L61").Value = 5675
$ws.Range("N61").Value = -6099
$ws.Range("H74").Value = 2472.8333
$ws.Range("I74").Value = 1402.0588
$ws.Range("J74").Value = 5073.2856
$ws.Range("K74").Value = 1402.0588
$ws.Range("L74").Value = 5073.2856
$ws.Range("M74").Value = -528.0588
$ws.Range("N74").Value = -6821.2856
$ws.Range("H77").Value = 2472.8333
$ws.Range("I77").Value = 1402.0588
$ws.Range("J77").Value = 5073.2856
$ws.Range("K77").Value = 7010.294
$ws.Range("L77").Value = 25366.428
$ws.Range("M77").Value = -2642.294
$ws.Range("N77").Value = -34102.428
$ws.Range("H122").Value = 2636.7778
$ws.Range("I122").Value = 1818.8572
$ws.Range("J122").Value = 5499.5
$ws.Range("K122").Value = 5456.571599999999
$ws.Range("L122").Value = 16498.5
$ws.Range("M122").Value = -3006.571599999999
$ws.Range("N122").Value = -21398.5
$ws.Range("H136").Value = 3693.5386
$ws.Range("J136").Value = 5675
$ws.Range("L136").Value = 17025
$ws.Range("N136").Value = -22125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9104.200000000001
$ws.Range("I99").Value = 4454.75
$ws.Range("J99").Value = 27702
$ws.Range("K99").Value = 4454.75
$ws.Range("L99").Value = 27702
$ws.Range("M99").Value = -2956.75
$ws.Range("N99").Value = -30698
$ws.Range("H105").Value = 1492.6666
$ws.Range("I105").Value = 1480
$ws.Range("J105").Value = 1499
$ws.Range("K105").Value = 1480
$ws.Range("L105").Value = 1499
$ws.Range("M105").Value = 267
$ws.Range("N105").Value = -4993
$ws.Range("H107").Value = 2687.3572
$ws.Range("I107").Value = 2528.1304
$ws.Range("K107").Value = 2528.1304
$ws.Range("M107").Value = -608.1304

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3299.6
$ws.Range("I16").Value = 3299.6
$ws.Range("K16").Value = 3299.6
$ws.Range("M16").Value = -3012.6
$ws.Range("H31").Value = 3599.158
$ws.Range("I31").Value = 1508.5555
$ws.Range("J31").Value = 5480.7
$ws.Range("K31").Value = 1508.5555
$ws.Range("L31").Value = 5480.7
$ws.Range("M31").Value = -1213.5555
$ws.Range("N31").Value = -6070.7
$ws.Range("H34").Value = 3599.158
$ws.Range("I34").Value = 1508.5555
$ws.Range("J34").Value = 5480.7
$ws.Range("K34").Value = 1508.5555
$ws.Range("L34").Value = 5480.7
$ws.Range("M34").Value = -1306.5555
$ws.Range("N34").Value = -5884.7
$ws.Range("H58").Value = 7624.375
$ws.Range("I58").Value = 3749.5
$ws.Range("K58").Value = 3749.5
$ws.Range("M58").Value = -3546.5
$ws.Range("H113").Value = 3299.6
$ws.Range("I113").Value = 3299.6
$ws.Range("K113").Value = 3299.6
$ws.Range("M113").Value = -1129.6
$ws.Range("H122").Value = 1314.1154
$ws.Range("I122").Value = 1049.0454
$ws.Range("J122").Value = 2772
$ws.Range("K122").Value = 3147.1362
$ws.Range("L122").Value = 8316
$ws.Range("M122").Value = -697.1361999999999
$ws.Range("N122").Value = -13216
$ws.Range("H132").Value = 40819284
$ws.Range("I132").Value = 48782810
$ws.Range("K132").Value = 146348430
$ws.Range("M132").Value = -146345900
$ws.Range("H134").Value = 2012.9642
$ws.Range("I134").Value = 1815.04
$ws.Range("K134").Value = 5445.12
$ws.Range("M134").Value = -2910.12
$ws.Range("H136").Value = 7624.375
$ws.Range("I136").Value = 3749.5
$ws.Range("K136").Value = 11248.5
$ws.Range("M136").Value = -8698.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.3125
$ws.Range("J2").Value = 60.11111
$ws.Range("L2").Value = 360.66666
$ws.Range("N2").Value = -586.66666
$ws.Range("H5").Value = 587.1667
$ws.Range("I5").Value = 587.1667
$ws.Range("K5").Value = 1761.5001
$ws.Range("M5").Value = -1649.5001
$ws.Range("H68").Value = 9257.375
$ws.Range("I68").Value = 2799.8
$ws.Range("J68").Value = 12192.637
$ws.Range("K68").Value = 8399.400000000001
$ws.Range("L68").Value = 36577.911
$ws.Range("M68").Value = -7588.400000000001
$ws.Range("N68").Value = -38199.911
$ws.Range("H71").Value = 9257.375
$ws.Range("I71").Value = 2799.8
$ws.Range("J71").Value = 12192.637
$ws.Range("K71").Value = 25198.2
$ws.Range("L71").Value = 109733.733
$ws.Range("M71").Value = -21142.2
$ws.Range("N71").Value = -117845.733
$ws.Range("H80").Value = 3966.5
$ws.Range("I80").Value = 3600
$ws.Range("J80").Value = 4149.75
$ws.Range("K80").Value = 10800
$ws.Range("L80").Value = 12449.25
$ws.Range("M80").Value = -9864
$ws.Range("N80").Value = -14321.25
$ws.Range("H83").Value = 3966.5
$ws.Range("I83").Value = 3600
$ws.Range("J83").Value = 4149.75
$ws.Range("K83").Value = 32400
$ws.Range("L83").Value = 37347.75
$ws.Range("M83").Value = -27720
$ws.Range("N83").Value = -46707.75
$ws.Range("H117").Value = 4591.357
$ws.Range("J117").Value = 5064.9165
$ws.Range("L117").Value = 15194.7495
$ws.Range("N117").Value = -22078.7495
$ws.Range("H122").Value = 16666894
$ws.Range("J122").Value = 28571640
$ws.Range("L122").Value = 257144760
$ws.Range("N122").Value = -257149660
$ws.Range("H129").Value = 1483.4
$ws.Range("I129").Value = 976.2857
$ws.Range("J129").Value = 2666.6667
$ws.Range("K129").Value = 2928.8571
$ws.Range("L129").Value = 8000.000100000001
$ws.Range("M129").Value = 2071.1429
$ws.Range("N129").Value = -18000.0001
$ws.Range("H135").Value = 587.1667
$ws.Range("I135").Value = 587.1667
$ws.Range("K135").Value = 5284.5003
$ws.Range("M135").Value = -2749.5003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6539.7
$ws.Range("J7").Value = 14532.667
$ws.Range("L7").Value = 14532.667
$ws.Range("N7").Value = -14756.667
$ws.Range("H16").Value = 6892.4375
$ws.Range("I16").Value = 11876
$ws.Range("J16").Value = 485
$ws.Range("K16").Value = 11876
$ws.Range("L16").Value = 485
$ws.Range("M16").Value = -11706
$ws.Range("N16").Value = -825
$ws.Range("H22").Value = 647.8333
$ws.Range("I22").Value = 586
$ws.Range("J22").Value = 833.3333
$ws.Range("K22").Value = 586
$ws.Range("L22").Value = 833.3333
$ws.Range("M22").Value = -291
$ws.Range("N22").Value = -1423.3333
$ws.Range("H27").Value = 647.8333
$ws.Range("I27").Value = 586
$ws.Range("J27").Value = 833.3333
$ws.Range("K27").Value = 586
$ws.Range("L27").Value = 833.3333
$ws.Range("M27").Value = -479
$ws.Range("N27").Value = -1047.3333
$ws.Range("H46").Value = 519.75
$ws.Range("J46").Value = 485
$ws.Range("L46").Value = 485
$ws.Range("N46").Value = -861
$ws.Range("H55").Value = 205.95
$ws.Range("I55").Value = 193.9
$ws.Range("J55").Value = 218
$ws.Range("K55").Value = 193.9
$ws.Range("L55").Value = 218
$ws.Range("M55").Value = -20.90000000000001
$ws.Range("N55").Value = -564
$ws.Range("H122").Value = 5571.357
$ws.Range("I122").Value = 4999.8887
$ws.Range("J122").Value = 6600
$ws.Range("K122").Value = 14999.6661
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -12549.6661
$ws.Range("N122").Value = -24700
$ws.Range("H126").Value = 6539.7
$ws.Range("J126").Value = 14532.667
$ws.Range("L126").Value = 43598.001
$ws.Range("N126").Value = -48538.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12727
$ws.Range("I81").Value = 13285.429
$ws.Range("J81").Value = 11749.75
$ws.Range("K81").Value = 26570.858
$ws.Range("L81").Value = 23499.5
$ws.Range("M81").Value = -25509.858
$ws.Range("N81").Value = -25621.5
$ws.Range("H84").Value = 12727
$ws.Range("I84").Value = 13285.429
$ws.Range("J84").Value = 11749.75
$ws.Range("K84").Value = 132854.29
$ws.Range("L84").Value = 117497.5
$ws.Range("M84").Value = -127550.29
$ws.Range("N84").Value = -128105.5
$ws.Range("H107").Value = 2472.0417
$ws.Range("I107").Value = 1340.9445
$ws.Range("K107").Value = 4022.8335
$ws.Range("M107").Value = -2102.8335
$ws.Range("H136").Value = 15878399
$ws.Range("I136").Value = 19613204
$ws.Range("J136").Value = 5474
$ws.Range("K136").Value = 58839612
$ws.Range("L136").Value = 16422
$ws.Range("M136").Value = -58837062

Write-Output "Applied 241 cell updates across 7 sheets"